$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("service")

# Add new row 10: Schema registry (set in this order so new shared-string
# entries land in the same order as the target workbook)
$ws.Range("A10").Value = "Schema registry"
$ws.Range("C10").Value = "linux-083:50014"

# Update Kafka connect address (row 5, column C) - shorten hostname
$ws.Range("C5").Value = "linux-084:50001"

# Update view: remove topLeftCell freeze by resetting scroll, change selection to C5
$ws.Range("C5").Select()
